# Updating filtered feeds from workflow
# Appends two new rows to the "Filtered Feeds" sheet for the
# "Qiagen, Incyte Partner to Develop Companion Diagnostics for CALR-Mutant,
# Other Blood Cancers" article (one row per source link).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$title = "Qiagen, Incyte Partner to Develop Companion Diagnostics for CALR-Mutant, Other Blood Cancers"
$keyword = "companion diagnostics"

$link1 = "https://www.genomeweb.com/cancer/qiagen-incyte-partner-develop-companion-diagnostics-calr-mutant-other-blood-cancers"
$link2 = "https://www.360dx.com/cancer/qiagen-incyte-partner-develop-companion-diagnostics-calr-mutant-other-blood-cancers"

# Row 23
$ws.Cells.Item(23, 1).Value = $link1
$ws.Cells.Item(23, 2).Value = $keyword
$ws.Cells.Item(23, 3).Value = $title

# Row 24
$ws.Cells.Item(24, 1).Value = $link2
$ws.Cells.Item(24, 2).Value = $keyword
$ws.Cells.Item(24, 3).Value = $title

# Wire up the hyperlinks on column A, matching the style used by the
# existing rows (the "Hyperlink" cell style).
$ws.Hyperlinks.Add($ws.Cells.Item(23, 1), $link1)
$ws.Cells.Item(23, 1).Style = $ws.Cells.Item(22, 1).Style

$ws.Hyperlinks.Add($ws.Cells.Item(24, 1), $link2)
$ws.Cells.Item(24, 1).Style = $ws.Cells.Item(22, 1).Style
